$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D:E (rows 2-51) to Text format so numeric-looking values
# (e.g. "6.14", "0.999") are stored as text, matching the original inlineStr cells,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '57.860.31'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '3.136.87'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '526.29'
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").Value = '141.73'
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.140.02'
$ws.Range("E8").Value = '  +2.44%  '
$ws.Range("E9").Value = '  +0.61%  '
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("E11").Value = '  +1.93%  '
$ws.Range("E12").Value = '  +3.28%  '
$ws.Range("D13").Value = '3.669.37'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("D14").Value = '0.132'
$ws.Range("E14").Value = '  +1.66%  '
$ws.Range("D15").Value = '26.45'
$ws.Range("E15").Value = '  +3.48%  '
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = '57.945.23'
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").Value = '3.125.03'
$ws.Range("E18").Value = '  +1.86%  '
$ws.Range("D19").Value = '6.14'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '12.96'
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("D21").Value = '8.13'
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").Value = '337.11'
$ws.Range("E22").Value = '  +0.99%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '0.513'
$ws.Range("E24").Value = '  +2.55%  '
$ws.Range("D25").Value = '66.79'
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").Value = '0.0₃0933'
$ws.Range("E28").Value = '  +2.85%  '
$ws.Range("D29").Value = '6.66'
$ws.Range("E29").Value = '  +4.66%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '7.24'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("E32").Value = '  +2.93%  '
$ws.Range("E33").Value = '  +2.22%  '
$ws.Range("D34").Value = '21.07'
$ws.Range("E34").Value = '  +1.15%  '
$ws.Range("D35").Value = '4.69'
$ws.Range("E35").Value = '  +5.10%  '
$ws.Range("D36").Value = '154.45'
$ws.Range("E36").Value = '  -0.48%  '
$ws.Range("D37").Value = '6.13'
$ws.Range("E37").Value = '  +3.57%  '
$ws.Range("D38").Value = '27.43'
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("E39").Value = '  +4.10%  '
$ws.Range("D40").Value = '0.0668'
$ws.Range("E40").Value = '  -0.89%  '
$ws.Range("D41").Value = '3.171.16'
$ws.Range("E41").Value = '  +2.14%  '
$ws.Range("D42").Value = '0.691'
$ws.Range("E42").Value = '  +5.38%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '3.92'
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '37.06'
$ws.Range("E44").Value = '  +0.37%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '1.51'
$ws.Range("E45").Value = '  +9.98%  '
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = '2.304.48'
$ws.Range("E47").Value = '  +2.01%  '
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("D49").Value = '0.994'
$ws.Range("E49").Value = '  +7.02%  '
$ws.Range("D50").Value = '20.97'
$ws.Range("E50").Value = '  +2.29%  '
$ws.Range("E51").Value = '  +2.96%  '

# Restore default (Normal) cell style so no stray number-format styling remains
$ws.Range("D2:E51").Style = "Normal"
